$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 296, pushing the existing data (old rows
# 296..379) down to 297..380.
$ws.Rows("296:296").Insert()

# Populate the newly inserted row 296 with the new weekly price record.
$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = 44876
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 100112001
$ws.Range("G296").Value = "Berenjena"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 65
$ws.Range("K296").Value = 13000
$ws.Range("L296").Value = 14000
$ws.Range("M296").Value = 13538
$ws.Range("N296").Value = "`$/caja 60 unidades"
$ws.Range("O296").Value = "Región de Arica y Parinacota"
$ws.Range("P296").Value = 226
$ws.Range("Q296").Value = 60
$ws.Range("R296").Value = "Hortaliza"
